# job_history.xlsx -- "fixed path in trajectory python code"
#
# A new trajectory entry ("5 batches", L/M = "-27.4 & -10.3" / "22.2 & 13.0")
# was logged ahead of the existing "filter large eigen vec..." block, and a
# few blank spacer rows were inserted along with it. Net effect: 4 brand new
# rows are inserted at row 53 (old rows 54-64 shift down to 58-68); the
# (previously blank) row 53 is rewritten with the new trajectory data, and
# the newly inserted rows 54-57 become the new spacer/partial rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 4 new rows at row 54, pushing the old rows 54..64 down to 58..68
#    and leaving the (still blank) old row 53 where it is.
$ws.Rows("54:57").Insert()

# 2) Row 53 becomes a populated data row (was fully blank before).
$ws.Range("A53").Value = "ukb51139_subset.csv"
$ws.Range("B53").Value = "28012 x 1081"
$ws.Range("C53").Value = "all"
$ws.Range("D53").Value = "no events"
$ws.Range("E53").Value = "> 160/100"
$ws.Range("F53").Value = "zscore"
$ws.Range("G53").Value = "median"
$ws.Range("H53").Value = "none"
$ws.Range("I53").Value = 25
$ws.Range("K53").Value = "N/A"
$ws.Range("L53").Value = "-27.4 & -10.3"
$ws.Range("M53").Value = "22.2 & 13.0"
$ws.Range("N53").Value = "N/A"
$ws.Range("O53").Value = "N/A"
$ws.Range("P53").Value = "5 batches"
$ws.Rows("53:53").RowHeight = 18.75

# 3) New row 54: same run-config columns A-H, I and the N/A markers, but no
#    L/M/P notes yet (those stay blank).
$ws.Range("A54").Value = "ukb51139_subset.csv"
$ws.Range("B54").Value = "28012 x 1081"
$ws.Range("C54").Value = "all"
$ws.Range("D54").Value = "no events"
$ws.Range("E54").Value = "> 160/100"
$ws.Range("F54").Value = "zscore"
$ws.Range("G54").Value = "median"
$ws.Range("H54").Value = "none"
$ws.Range("I54").Value = 25
$ws.Range("K54").Value = "N/A"
$ws.Range("N54").Value = "N/A"
$ws.Range("O54").Value = "N/A"
$ws.Rows("54:54").RowHeight = 18.75

# 4) New rows 55 and 56 stay fully blank (spacer rows), just with a slightly
#    different row height than the rest of the table and right-aligned,
#    unbordered number formats on the numeric columns.
foreach ($r in 55, 56) {
    $ws.Range("I$r").NumberFormat = "#,##0"
    $ws.Range("I$r").Borders.LineStyle = -4142
    $ws.Range("K$r").NumberFormat = "#,##0"
    $ws.Range("K$r").Borders.LineStyle = -4142
    $ws.Range("N$r").NumberFormat = "#,##0"
    $ws.Range("N$r").Borders.LineStyle = -4142
    $ws.Range("O$r").NumberFormat = "#,##0.00"
    $ws.Range("O$r").Borders.LineStyle = -4142
    $ws.Range("P$r").Borders.LineStyle = -4142
    $ws.Rows("${r}:${r}").RowHeight = 18.75
}

# 5) New row 57 stays fully blank, with the regular table row height/style
#    (this one already matches the default from the insert, nothing to do).

"done"
